# Apply the "Updated cryptos list" data refresh described by the commit diff.
# Each row is the (cell, new value) pair taken from the unified OOXML diff.
# Numeric-looking strings are written with a leading apostrophe so Excel
# keeps them as literal text (matching the original inline-string cells)
# instead of silently coercing them to floating point numbers, which would
# corrupt values like "1.000", "7.100", or "30.773.73".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.773.73'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.890.82'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range("D5").Value = "'249.46"
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range("D6").Value = "'0.9994"
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range("D8").Value = "'0.2936"
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range("D9").Value = "'0.06538"
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range("D10").Value = "'22.07"
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range("D11").Value = "'0.07759"
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range("D12").Value = "'97.34"
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.889.77'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = "'0.7371"
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range("D15").Value = "'5.248"
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range("D16").Value = "'283.99"
$ws.Range('E16').Value = '  +3.24%  '
$ws.Range('D17').Value = '30.818.77'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range("D19").Value = "'0.000007585"
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range("D20").Value = "'0.9998"
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = '2.136.71'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range("D22").Value = "'5.343"
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range("D23").Value = "'0.9997"
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range("D24").Value = "'6.253"
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range("D25").Value = "'9.257"
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range("D26").Value = "'164.28"
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  -1.08%  '
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range("D30").Value = "'0.09752"
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range("D31").Value = "'1.501"
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range("D32").Value = "'4.310"
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range("D33").Value = "'4.193"
$ws.Range('E33').Value = '  +1.45%  '
$ws.Range("D34").Value = "'0.04884"
$ws.Range('E34').Value = '  +1.62%  '
$ws.Range('E35').Value = '  -0.34%  '
$ws.Range("D36").Value = "'0.7009"
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range("D37").Value = "'2.720"
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range("D38").Value = "'0.01916"
$ws.Range('E38').Value = '  +2.35%  '
$ws.Range('E39').Value = '  +2.07%  '
$ws.Range("D40").Value = "'6.357"
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range("D41").Value = "'76.14"
$ws.Range('E41').Value = '  +6.50%  '
$ws.Range("D42").Value = "'2.030"
$ws.Range('E42').Value = '  +2.19%  '
$ws.Range("D43").Value = "'0.4267"
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = "'0.9994"
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = "'0.8385"
$ws.Range('E45').Value = '  -0.36%  '
$ws.Range("D46").Value = "'102.09"
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range("D47").Value = "'9.442"
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range("D48").Value = "'7.096"
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range("D49").Value = "'35.85"
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range("D50").Value = "'927.61"
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range("D51").Value = "'0.05774"
$ws.Range('E51').Value = '  +2.24%  '
